$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45 (shifts existing rows 45-48 down to 46-49)
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new weekly price record
$ws.Cells.Item(45, 1).Value = 4
$ws.Cells.Item(45, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(45, 3).Value = "Los Lagos"
$ws.Cells.Item(45, 4).Value = 45223
$ws.Cells.Item(45, 5).Value = 10
$ws.Cells.Item(45, 6).Value = 100112013
$ws.Cells.Item(45, 7).Value = "Alcachofa"
$ws.Cells.Item(45, 8).Value = "Española"
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 250
$ws.Cells.Item(45, 11).Value = 13500
$ws.Cells.Item(45, 12).Value = 14000
$ws.Cells.Item(45, 13).Value = 13700
$ws.Cells.Item(45, 14).Value = "`$/caja 30 unidades"
$ws.Cells.Item(45, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(45, 16).Value = 457
$ws.Cells.Item(45, 17).Value = 30
$ws.Cells.Item(45, 18).Value = "Hortaliza"
